$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a literal TEXT value into a cell while preserving / forcing a
# specific cell style (copied from a reference cell that already carries the
# desired style in this sheet). Setting NumberFormatLocal to "@" before the
# value assignment keeps numeric-looking strings (like "0") stored as text
# instead of being coerced to a number; the subsequent PasteSpecial (formats
# only) from the reference cell then restores the exact visual style/index
# expected for this column.
# ---------------------------------------------------------------------------
function Set-TextCell($ws, [string]$addr, [string]$text, [string]$styleRef) {
    $dst = $ws.Range($addr)
    $dst.NumberFormatLocal = "@"
    $dst.Value = $text
    $src = $ws.Range($styleRef)
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

function Set-NumCell($ws, [string]$addr, $num, [string]$styleRef) {
    $dst = $ws.Range($addr)
    $dst.Value = $num
    if ($styleRef) {
        $src = $ws.Range($styleRef)
        $src.Copy()
        $dst.PasteSpecial(-4122)  # xlPasteFormats
    }
}

# Reference cells for each cell style used below (row 15 values untouched by
# the diff so they're safe, stable style donors).
$styleText = "C15"   # s=14 (text, e.g. "0"/"***.*")
$styleNum15 = "I15"  # s=15 (plain integer count)
$styleNum16 = "L15"  # s=16 (percentage-style number)

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Range("N15").Value = -50

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
Set-TextCell $ws "C16" "0" $styleText
# D16 and E16 remain "0" / "***.*" text (unchanged by the diff; left as-is)
Set-NumCell $ws "F16" 5 $null
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = -42.307692307692

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
Set-NumCell $ws "D17" 1 $styleNum15
Set-NumCell $ws "E17" 0 $styleNum16
$ws.Range("F17").Value = 2
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 3
$ws.Range("J17").Value = 3
$ws.Range("L17").Value = 50
Set-NumCell $ws "M17" 200 $styleNum16
$ws.Range("N17").Value = -72.727272727272

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
Set-TextCell $ws "D19" "0" $styleText
Set-TextCell $ws "E19" "***.*" $styleText
$ws.Range("F19").Value = 6
$ws.Range("H19").Value = 500
$ws.Range("M19").Value = 0

# ---------------------------------------------------------------------------
# Row 21 (TOTAL) - values only, styles unchanged
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 13
$ws.Range("H21").Value = 225
$ws.Range("I21").Value = 28
$ws.Range("J21").Value = 11
$ws.Range("K21").Value = 154.545454545455
$ws.Range("L21").Value = 154.545454545455
$ws.Range("M21").Value = 75
$ws.Range("N21").Value = -55.555555555555

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = -50
$ws.Range("M24").Value = -50

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
Set-TextCell $ws "C26" "0" $styleText
Set-TextCell $ws "D26" "0" $styleText
Set-TextCell $ws "E26" "***.*" $styleText
$ws.Range("M26").Value = -50

# ---------------------------------------------------------------------------
# Row 28
# ---------------------------------------------------------------------------
$ws.Range("D28").Value = 2
Set-TextCell $ws "F28" "0" $styleText
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -100
$ws.Range("J28").Value = 5
$ws.Range("K28").Value = -80

# ---------------------------------------------------------------------------
# Shared-string rich-text edits (partial run replacement via Characters so
# the remaining runs / their formatting stay intact).
# ---------------------------------------------------------------------------

# "Volume 31   Number  13" -> "...14"  (A8)
$a8 = $ws.Range("A8")
$a8v = $a8.Value2
$start = $a8v.IndexOf("13") + 1
$len = "13".Length
$a8.Characters($start, $len).Text = "14"

# "Report Covering the Week  3/25/2024  Through  3/31/2024"
#   -> "...4/1/2024  Through  4/7/2024"  (C9)
$c9 = $ws.Range("C9")
$c9v = $c9.Value2
$startA = $c9v.IndexOf("3/25/2024") + 1
$lenA = "3/25/2024".Length
$c9.Characters($startA, $lenA).Text = "4/1/2024"

$c9v2 = $c9.Value2
$startB = $c9v2.IndexOf("3/31/2024") + 1
$lenB = "3/31/2024".Length
$c9.Characters($startB, $lenB).Text = "4/7/2024"
